$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the hidden/encrypted value in column B (rows 2-11) with the new
# encrypted string value. Since every cell referencing the old shared
# string gets overwritten, the old shared string entry is dropped and the
# new one is appended at the end of the shared strings table.
$newValue = "U2FsdGVkX1+s1oRPt0B7OcjpVP5f3IBhIA53DP6hmaIufy/vOcfu72zH6t1I7wrJH9m3kgsE3/e9RQkIhB+/5g=="

for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 2).Value = $newValue
}

# Update the active selection on the sheet from B16 to B12.
$ws.Range("B12").Select()
